$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "MEC-3B-Mec. Manut. Equip. Ind."
$ws.Range("D2").Value = "-"

$ws.Range("C3").Value = "MEC-3B-M. S. Ar Cond."
$ws.Range("D3").Value = "-"

$ws.Range("B4").Value = "MEC-3B-Mec. Manut. Equip. Ind."
$ws.Range("D4").Value = "-"

$ws.Range("B6").Value = "MEC-3B-Mec. Manut. Equip. Ind."
$ws.Range("D6").Value = "-"
$ws.Range("E6").Value = "-"

$ws.Range("B7").Value = "MEC-3B-Mec. Manut. Equip. Ind."
$ws.Range("D7").Value = "MEC-3B-M. S. Ar Cond."
$ws.Range("E7").Value = "-"

$ws.Range("B8").Value = "MEC-3B-M. S. Ar Cond."
$ws.Range("D8").Value = "MEC-3B-M. S. Ar Cond."
$ws.Range("E8").Value = "-"
